$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 13:35"

# España (row 5)
$ws.Range("B5").Value = 272646
$ws.Range("C5").Value = 1551
$ws.Range("D5").Value = 186480
$ws.Range("E5").Value = 58845
$ws.Range("F5").Value = 1376
$ws.Range("G5").Value = 217
$ws.Range("H5").Value = 27321

# Republica de Macedonia (row 87)
$ws.Range("B87").Value = 1723
$ws.Range("C87").Value = 29
$ws.Range("D87").Value = 1235
$ws.Range("E87").Value = 393

# Rows 123-125 re-sort: Malta overtakes Jamaica & Tanzania
$ws.Range("A123").Value = "Malta"
$ws.Range("B123").Value = 522
$ws.Range("C123").Value = 14
$ws.Range("D123").Value = 443
$ws.Range("E123").Value = 73
$ws.Range("F123").Value = 1
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 6

$ws.Range("A124").Value = "Jamaica"
$ws.Range("B124").Value = 509
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 113
$ws.Range("E124").Value = 387
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 9

$ws.Range("A125").Value = "Tanzania"
$ws.Range("B125").Value = 509
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 183
$ws.Range("E125").Value = 305
$ws.Range("F125").Value = 7
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 21
